$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.018677635229236
$ws.Cells.Item(2, 4).Value2 = 1.054350528882185
$ws.Cells.Item(2, 5).Value2 = 1.019891571743906
$ws.Cells.Item(2, 6).Value2 = 1.056441662908779
$ws.Cells.Item(2, 9).Value2 = 1.04368904623742
$ws.Cells.Item(2, 10).Value2 = 1.023884710911551
$ws.Cells.Item(2, 11).Value2 = 1.057093980672994
$ws.Cells.Item(2, 12).Value2 = 1.022733475487947
$ws.Cells.Item(2, 13).Value2 = 1.059179364208686
$ws.Cells.Item(2, 14).Value2 = 1.011967439458822

$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.020208398964408
$ws.Cells.Item(3, 4).Value2 = 1.055181212700848
$ws.Cells.Item(3, 5).Value2 = 1.02120501077031
$ws.Cells.Item(3, 6).Value2 = 1.057546213006259
$ws.Cells.Item(3, 9).Value2 = 1.043942093781302
$ws.Cells.Item(3, 10).Value2 = 1.025048857346466
$ws.Cells.Item(3, 11).Value2 = 1.057737845677014
$ws.Cells.Item(3, 12).Value2 = 1.023851643664048
$ws.Cells.Item(3, 13).Value2 = 1.060096813535946
$ws.Cells.Item(3, 14).Value2 = 1.012358264058563

$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.021196670310999
$ws.Cells.Item(4, 4).Value2 = 1.055715125493566
$ws.Cells.Item(4, 5).Value2 = 1.022053318904646
$ws.Cells.Item(4, 6).Value2 = 1.05825676394047
$ws.Cells.Item(4, 9).Value2 = 1.044102402216575
$ws.Cells.Item(4, 10).Value2 = 1.025799604177475
$ws.Cells.Item(4, 11).Value2 = 1.058150234197895
$ws.Cells.Item(4, 12).Value2 = 1.024573072636481
$ws.Cells.Item(4, 13).Value2 = 1.060685716629813
$ws.Cells.Item(4, 14).Value2 = 1.012610124411189

$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.021611613527143
$ws.Cells.Item(5, 4).Value2 = 1.055938723062466
$ws.Cells.Item(5, 5).Value2 = 1.022409577277587
$ws.Cells.Item(5, 6).Value2 = 1.058554485113095
$ws.Cells.Item(5, 9).Value2 = 1.044168975126693
$ws.Cells.Item(5, 10).Value2 = 1.026114618021545
$ws.Cells.Item(5, 11).Value2 = 1.058322589645549
$ws.Cells.Item(5, 12).Value2 = 1.024875864519336
$ws.Cells.Item(5, 13).Value2 = 1.060932157466119
$ws.Cells.Item(5, 14).Value2 = 1.012715762238394

$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.021681253731393
$ws.Cells.Item(6, 4).Value2 = 1.055976215720745
$ws.Cells.Item(6, 5).Value2 = 1.022469373098571
$ws.Cells.Item(6, 6).Value2 = 1.058604415581268
$ws.Cells.Item(6, 9).Value2 = 1.044180104893613
$ws.Cells.Item(6, 10).Value2 = 1.026167475226575
$ws.Cells.Item(6, 11).Value2 = 1.058351469532761
$ws.Cells.Item(6, 12).Value2 = 1.024926675646599
$ws.Cells.Item(6, 13).Value2 = 1.060973469541542
$ws.Cells.Item(6, 14).Value2 = 1.012733485032532

$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.021202216858586
$ws.Cells.Item(7, 4).Value2 = 1.055718116590602
$ws.Cells.Item(7, 5).Value2 = 1.022058080693682
$ws.Cells.Item(7, 6).Value2 = 1.058260746009314
$ws.Cells.Item(7, 9).Value2 = 1.044103294991682
$ws.Cells.Item(7, 10).Value2 = 1.025803815758126
$ws.Cells.Item(7, 11).Value2 = 1.058152541197755
$ws.Cells.Item(7, 12).Value2 = 1.024577120501429
$ws.Cells.Item(7, 13).Value2 = 1.060689014037341
$ws.Cells.Item(7, 14).Value2 = 1.012611536906191

$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.019195432275664
$ws.Cells.Item(8, 4).Value2 = 1.054632007604759
$ws.Cells.Item(8, 5).Value2 = 1.020335784688701
$ws.Cells.Item(8, 6).Value2 = 1.056815815469327
$ws.Cells.Item(8, 9).Value2 = 1.043775275385672
$ws.Cells.Item(8, 10).Value2 = 1.024278668197809
$ws.Cells.Item(8, 11).Value2 = 1.057312455560118
$ws.Cells.Item(8, 12).Value2 = 1.02311180425784
$ws.Cells.Item(8, 13).Value2 = 1.059490405237387
$ws.Cells.Item(8, 14).Value2 = 1.012099734887991

$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.015641693769673
$ws.Cells.Item(9, 4).Value2 = 1.052690544564786
$ws.Cells.Item(9, 5).Value2 = 1.017288517935886
$ws.Cells.Item(9, 6).Value2 = 1.054237603721856
$ws.Cells.Item(9, 9).Value2 = 1.04317096752418
$ws.Cells.Item(9, 10).Value2 = 1.021571464362425
$ws.Cells.Item(9, 11).Value2 = 1.055799619961755
$ws.Cells.Item(9, 12).Value2 = 1.020513368018558
$ws.Cells.Item(9, 13).Value2 = 1.057341813743756
$ws.Cells.Item(9, 14).Value2 = 1.011189893571253

$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.013260172105644
$ws.Cells.Item(10, 4).Value2 = 1.051377576428709
$ws.Cells.Item(10, 5).Value2 = 1.015248270099376
$ws.Cells.Item(10, 6).Value2 = 1.052497023205555
$ws.Cells.Item(10, 9).Value2 = 1.042750376236787
$ws.Cells.Item(10, 10).Value2 = 1.019753000881598
$ws.Cells.Item(10, 11).Value2 = 1.054769117685309
$ws.Cells.Item(10, 12).Value2 = 1.018769688371083
$ws.Cells.Item(10, 13).Value2 = 1.055884697129749
$ws.Cells.Item(10, 14).Value2 = 1.010577828407062

$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.012225884343118
$ws.Cells.Item(11, 4).Value2 = 1.050804594875981
$ws.Cells.Item(11, 5).Value2 = 1.014362651712469
$ws.Cells.Item(11, 6).Value2 = 1.051738116279871
$ws.Cells.Item(11, 9).Value2 = 1.042564041977656
$ws.Cells.Item(11, 10).Value2 = 1.018962251493866
$ws.Cells.Item(11, 11).Value2 = 1.054317669259246
$ws.Cells.Item(11, 12).Value2 = 1.018011866647532
$ws.Cells.Item(11, 13).Value2 = 1.055247840066288
$ws.Cells.Item(11, 14).Value2 = 1.010311460493468

$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.01184123072662
$ws.Cells.Item(12, 4).Value2 = 1.05059109175183
$ws.Cells.Item(12, 5).Value2 = 1.014033358008769
$ws.Cells.Item(12, 6).Value2 = 1.051455434767651
$ws.Cells.Item(12, 9).Value2 = 1.042494194692604
$ws.Cells.Item(12, 10).Value2 = 1.018668021460643
$ws.Cells.Item(12, 11).Value2 = 1.054149192441886
$ws.Cells.Item(12, 12).Value2 = 1.01772994999279
$ws.Cells.Item(12, 13).Value2 = 1.055010389914801
$ws.Cells.Item(12, 14).Value2 = 1.010212315558985

$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.011923761825652
$ws.Cells.Item(13, 4).Value2 = 1.050636919360129
$ws.Cells.Item(13, 5).Value2 = 1.014104007949475
$ws.Cells.Item(13, 6).Value2 = 1.05151610672585
$ws.Cells.Item(13, 9).Value2 = 1.042509205900958
$ws.Cells.Item(13, 10).Value2 = 1.018731158037372
$ws.Cells.Item(13, 11).Value2 = 1.054185367024376
$ws.Cells.Item(13, 12).Value2 = 1.017790441561742
$ws.Cells.Item(13, 13).Value2 = 1.055061364262357
$ws.Cells.Item(13, 14).Value2 = 1.010233591761593

$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.01219409847243
$ws.Cells.Item(14, 4).Value2 = 1.050786960370259
$ws.Cells.Item(14, 5).Value2 = 1.014335439087629
$ws.Cells.Item(14, 6).Value2 = 1.051714765875813
$ws.Cells.Item(14, 9).Value2 = 1.04255828132848
$ws.Cells.Item(14, 10).Value2 = 1.018937940810072
$ws.Cells.Item(14, 11).Value2 = 1.054303759017957
$ws.Cells.Item(14, 12).Value2 = 1.017988572103469
$ws.Cells.Item(14, 13).Value2 = 1.05522823060657
$ws.Cells.Item(14, 14).Value2 = 1.010303269320552

$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.012360598770793
$ws.Cells.Item(15, 4).Value2 = 1.050879316428965
$ws.Cells.Item(15, 5).Value2 = 1.014477986724946
$ws.Cells.Item(15, 6).Value2 = 1.051837061546009
$ws.Cells.Item(15, 9).Value2 = 1.042588434214362
$ws.Cells.Item(15, 10).Value2 = 1.019065278604144
$ws.Cells.Item(15, 11).Value2 = 1.054376599639031
$ws.Cells.Item(15, 12).Value2 = 1.018110589932473
$ws.Cells.Item(15, 13).Value2 = 1.05533092399601
$ws.Cells.Item(15, 14).Value2 = 1.010346172844836

$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.013328748661838
$ws.Cells.Item(16, 4).Value2 = 1.051415509194248
$ws.Cells.Item(16, 5).Value2 = 1.015306999077334
$ws.Cells.Item(16, 6).Value2 = 1.052547278884577
$ws.Cells.Item(16, 9).Value2 = 1.04276265371895
$ws.Cells.Item(16, 10).Value2 = 1.019805409181646
$ws.Cells.Item(16, 11).Value2 = 1.054798968350857
$ws.Cells.Item(16, 12).Value2 = 1.018819922918669
$ws.Cells.Item(16, 13).Value2 = 1.055926838231541
$ws.Cells.Item(16, 14).Value2 = 1.01059547790821

$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.013935212746376
$ws.Cells.Item(17, 4).Value2 = 1.051750653538939
$ws.Cells.Item(17, 5).Value2 = 1.015826427691681
$ws.Cells.Item(17, 6).Value2 = 1.052991377459829
$ws.Cells.Item(17, 9).Value2 = 1.042870807643546
$ws.Cells.Item(17, 10).Value2 = 1.020268772781056
$ws.Cells.Item(17, 11).Value2 = 1.055062506097057
$ws.Cells.Item(17, 12).Value2 = 1.019264114528129
$ws.Cells.Item(17, 13).Value2 = 1.056299052644987
$ws.Cells.Item(17, 14).Value2 = 1.0107514997847

$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.014288657247293
$ws.Cells.Item(18, 4).Value2 = 1.05194570764312
$ws.Cells.Item(18, 5).Value2 = 1.01612919208351
$ws.Cells.Item(18, 6).Value2 = 1.053249908934067
$ws.Cells.Item(18, 9).Value2 = 1.042933485347768
$ws.Cells.Item(18, 10).Value2 = 1.020538722974653
$ws.Cells.Item(18, 11).Value2 = 1.05521571835705
$ws.Cells.Item(18, 12).Value2 = 1.019522934643296
$ws.Cells.Item(18, 13).Value2 = 1.056515588509914
$ws.Cells.Item(18, 14).Value2 = 1.010842375636863

$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.014409122867677
$ws.Cells.Item(19, 4).Value2 = 1.052012143204401
$ws.Cells.Item(19, 5).Value2 = 1.016232391630266
$ws.Cells.Item(19, 6).Value2 = 1.053337976190734
$ws.Cells.Item(19, 9).Value2 = 1.042954787882771
$ws.Cells.Item(19, 10).Value2 = 1.020630714654438
$ws.Cells.Item(19, 11).Value2 = 1.055267874217567
$ws.Cells.Item(19, 12).Value2 = 1.019611140161651
$ws.Cells.Item(19, 13).Value2 = 1.056589325002972
$ws.Cells.Item(19, 14).Value2 = 1.010873340144609

$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.013870175561524
$ws.Cells.Item(20, 4).Value2 = 1.051714740209708
$ws.Cells.Item(20, 5).Value2 = 1.015770719650955
$ws.Cells.Item(20, 6).Value2 = 1.052943782021644
$ws.Cells.Item(20, 9).Value2 = 1.042859245820771
$ws.Cells.Item(20, 10).Value2 = 1.020219091613278
$ws.Cells.Item(20, 11).Value2 = 1.055034283227118
$ws.Cells.Item(20, 12).Value2 = 1.019216484886617
$ws.Cells.Item(20, 13).Value2 = 1.056259176601838
$ws.Cells.Item(20, 14).Value2 = 1.010734773484078

$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.012114504219626
$ws.Cells.Item(21, 4).Value2 = 1.050742795611229
$ws.Cells.Item(21, 5).Value2 = 1.01426729771834
$ws.Cells.Item(21, 6).Value2 = 1.051656287546434
$ws.Cells.Item(21, 9).Value2 = 1.042543847361339
$ws.Cells.Item(21, 10).Value2 = 1.018877062597046
$ws.Cells.Item(21, 11).Value2 = 1.054268917320466
$ws.Cells.Item(21, 12).Value2 = 1.017930239472526
$ws.Cells.Item(21, 13).Value2 = 1.055179117315733
$ws.Cells.Item(21, 14).Value2 = 1.010282756670746

$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.011007901071026
$ws.Cells.Item(22, 4).Value2 = 1.05012780461182
$ws.Cells.Item(22, 5).Value2 = 1.01332009072181
$ws.Cells.Item(22, 6).Value2 = 1.050842217266792
$ws.Cells.Item(22, 9).Value2 = 1.042341872528793
$ws.Cells.Item(22, 10).Value2 = 1.0180303178821
$ws.Cells.Item(22, 11).Value2 = 1.053783136540591
$ws.Cells.Item(22, 12).Value2 = 1.017119046189612
$ws.Cells.Item(22, 13).Value2 = 1.054494871853689
$ws.Cells.Item(22, 14).Value2 = 1.009997373901181

$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.011594796168842
$ws.Cells.Item(23, 4).Value2 = 1.050454192664077
$ws.Cells.Item(23, 5).Value2 = 1.013822410373763
$ws.Cells.Item(23, 6).Value2 = 1.051274206291987
$ws.Cells.Item(23, 9).Value2 = 1.042449291596534
$ws.Cells.Item(23, 10).Value2 = 1.018479476443726
$ws.Cells.Item(23, 11).Value2 = 1.054041091648991
$ws.Cells.Item(23, 12).Value2 = 1.017549312751778
$ws.Cells.Item(23, 13).Value2 = 1.054858094776435
$ws.Cells.Item(23, 14).Value2 = 1.010148773665564

$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.013899563979791
$ws.Cells.Item(24, 4).Value2 = 1.051730969226191
$ws.Cells.Item(24, 5).Value2 = 1.015795892364253
$ws.Cells.Item(24, 6).Value2 = 1.052965289908325
$ws.Cells.Item(24, 9).Value2 = 1.042864471367545
$ws.Cells.Item(24, 10).Value2 = 1.020241541389463
$ws.Cells.Item(24, 11).Value2 = 1.05504703748753
$ws.Cells.Item(24, 12).Value2 = 1.0192380075038
$ws.Cells.Item(24, 13).Value2 = 1.056277196632529
$ws.Cells.Item(24, 14).Value2 = 1.010742331778228

$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.01656255652206
$ws.Cells.Item(25, 4).Value2 = 1.053195740202022
$ws.Cells.Item(25, 5).Value2 = 1.018077817752377
$ws.Cells.Item(25, 6).Value2 = 1.054907954339736
$ws.Cells.Item(25, 9).Value2 = 1.04317096752418
$ws.Cells.Item(25, 10).Value2 = 1.022273718352478
$ws.Cells.Item(25, 11).Value2 = 1.056194584778136
$ws.Cells.Item(25, 12).Value2 = 1.021187104627113
$ws.Cells.Item(25, 13).Value2 = 1.057901619248833
$ws.Cells.Item(25, 14).Value2 = 1.01142606925239
